$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G10" = 1.48
    "I10" = 5.75
    "Q10" = 11
    "R10" = 12
    "Y10" = 17
    "AA10" = 17
    "AF10" = 15
    "G11" = 2
    "H11" = 3.2
    "I11" = 3.6
    "J11" = 2.32
    "K11" = 1.47
    "L11" = 1.5
    "M11" = 2.25
    "N11" = 5.6
    "O11" = 8.25
    "P11" = 9.25
    "Q11" = 17
    "R11" = 19.5
    "S11" = 40
    "T11" = 7
    "U11" = 6.3
    "V11" = 20
    "W11" = 120
    "Y11" = 8
    "Z11" = 17.5
    "AA11" = 13.5
    "AB11" = 55
    "AC11" = 40
    "AD11" = 60
    "AG11" = 1.45
    "AH11" = 2.37
    "AI11" = 2.07
    "AJ11" = 1.6
    "G12" = 2.25
    "H12" = 3.1
    "I12" = 3.1
    "J12" = 2.35
    "K12" = 1.47
    "L12" = 1.5
    "M12" = 2.25
    "N12" = 5.9
    "O12" = 9.5
    "P12" = 9.75
    "Q12" = 22
    "R12" = 23
    "S12" = 45
    "T12" = 6.8
    "U12" = 6.1
    "V12" = 19
    "W12" = 120
    "Y12" = 7.2
    "Z12" = 14
    "AA12" = 12
    "AB12" = 40
    "AC12" = 35
    "AD12" = 55
    "AG12" = 1.47
    "AH12" = 2.32
    "AI12" = 2.05
    "AJ12" = 1.6
    "G15" = 1.22
    "H15" = 4.9
    "I15" = 11.25
    "J15" = 1.65
    "K15" = 2
    "N15" = 5.5
    "O15" = 4.75
    "P15" = 7.8
    "Q15" = 5.8
    "R15" = 9.5
    "T15" = 11.5
    "U15" = 9
    "V15" = 22
    "Y15" = 22
    "Z15" = 70
    "AA15" = 29
    "AB15" = 300
    "AD15" = 100
    "J17" = 2.15
    "K17" = 1.67
    "AF17" = 9
    "G19" = 1.4
    "H19" = 3.7
    "I19" = 8
    "J19" = 2.15
    "K19" = 1.67
    "N19" = 5.5
    "T19" = 8
    "Y19" = 17
    "AG19" = 1.36
    "AH19" = 3
    "AI19" = 2.38
    "AJ19" = 1.53
    "J20" = 2.25
    "K20" = 1.62
    "L20" = 1.44
    "M20" = 2.63
    "R20" = 15
    "S20" = 41
    "T20" = 7.5
    "V20" = 23
    "AC20" = 67
    "AG20" = 1.4
    "AH20" = 2.75
    "AI20" = 2.38
    "AJ20" = 1.53
    "G21" = 1.44
    "H21" = 3.5
    "I21" = 8
    "J21" = 2.25
    "K21" = 1.62
    "L21" = 1.44
    "M21" = 2.4
    "P21" = 9.5
    "Q21" = 9
    "S21" = 41
    "V21" = 23
    "Y21" = 17
    "Z21" = 41
    "AA21" = 26
    "AB21" = 101
    "AC21" = 67
    "AD21" = 81
    "AG21" = 1.39
    "AH21" = 2.57
    "AI21" = 2.32
    "AJ21" = 1.47
    "G23" = 2.4
    "I23" = 2.9
    "J23" = 2.7
    "K23" = 1.44
    "L23" = 1.57
    "M23" = 2.25
    "N23" = 6
    "O23" = 10
    "Q23" = 23
    "R23" = 26
    "T23" = 6
    "W23" = 81
    "AG23" = 1.53
    "AH23" = 2.38
    "AI23" = 2.2
    "AJ23" = 1.62
    "G24" = 1.62
    "H24" = 3.75
    "I24" = 4.5
    "J24" = 2.1
    "K24" = 1.7
    "O24" = 7
    "P24" = 8.5
    "Q24" = 12
    "T24" = 9.5
    "Y24" = 11
    "AA24" = 17
    "AD24" = 51
    "AI24" = 2.02
    "AJ24" = 1.62
    "L25" = 1.53
    "M25" = 2.2
    "AG25" = 1.5
    "AH25" = 2.25
    "AI25" = 2.1
    "AJ25" = 1.57
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
